$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header cells: "<field>_old" -> "<field>_FV2304", "<field>_new" -> "<field>_FV2310"
$headers = @(
  "Segmentname_FV2304","Segmentgruppe_FV2304","Segment_FV2304","Datenelement_FV2304","Segment ID_FV2304",
  "Code_FV2304","Qualifier_FV2304","Beschreibung_FV2304","Bedingungsausdruck_FV2304","Bedingung_FV2304",
  "diff",
  "Segmentname_FV2310","Segmentgruppe_FV2310","Segment_FV2310","Datenelement_FV2310","Segment ID_FV2310",
  "Code_FV2310","Qualifier_FV2310","Beschreibung_FV2310","Bedingungsausdruck_FV2310","Bedingung_FV2310"
)

for ($i = 0; $i -lt $headers.Count; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2. Turn the used range A1:U61 into an Excel Table ("Table1"), headers already in place.
$tableRange = $ws.Range("A1:U61")
$lo = $ws.ListObjects.Add(1, $tableRange, [System.Reflection.Missing]::Value, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# 3. Freeze the header row (split below row 1, freeze panes).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
